# Recalculated extrapolation results after removing sub-$5 price noise
# from the calibration input. Updates ABSM1_RN/M1_RN/CM2_RN/CMN3_RN/CMN4_RN
# (columns D-H) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 115162.4317729514
$ws.Range("E3").Value = 0.03011455478647912
$ws.Range("F3").Value = 0.1664106990635528
$ws.Range("G3").Value = -1.401513319294426
$ws.Range("H3").Value = 13.22155647111928

$ws.Range("D5").Value = 116948.9686272062
$ws.Range("E5").Value = 0.01675859769807713
$ws.Range("F5").Value = 0.1821901062908523
$ws.Range("G5").Value = -0.7533364570375949
$ws.Range("H5").Value = 8.277389691029969

$ws.Range("D7").Value = 118625.947238204
$ws.Range("E7").Value = -0.00348448750378783
$ws.Range("F7").Value = 0.2317446559831317
$ws.Range("G7").Value = -0.9180057570673861
$ws.Range("H7").Value = 6.10295549526394

$ws.Range("D8").Value = 119053.2255420037
$ws.Range("E8").Value = -0.0257718284828914
$ws.Range("F8").Value = 0.2052184675360208
$ws.Range("G8").Value = -0.9716877234413515
$ws.Range("H8").Value = 7.072731469214665

$ws.Range("D9").Value = 120569.8317311448
$ws.Range("E9").Value = -0.05735926553632831
$ws.Range("F9").Value = 0.3245563305064358
$ws.Range("G9").Value = -1.609524609409163
$ws.Range("H9").Value = 10.48917018802919

$ws.Range("D10").Value = 122073.3437649611
$ws.Range("E10").Value = -0.1001216335303524
$ws.Range("F10").Value = 0.4382550954376864
$ws.Range("G10").Value = -1.893801721007187
$ws.Range("H10").Value = 9.623707528326589

$ws.Range("D11").Value = 124048.0995439282
$ws.Range("E11").Value = -0.1760418082382996
$ws.Range("F11").Value = 0.7560233799149142
$ws.Range("G11").Value = -2.553318547547716
$ws.Range("H11").Value = 12.26682259271791

$ws.Range("D13").Value = 114499.0498290738
$ws.Range("E13").Value = 0.08223135338759768
$ws.Range("F13").Value = 0.1525343821248598
$ws.Range("G13").Value = -0.7654399439376441
$ws.Range("H13").Value = 6.764653780905348

$ws.Range("D14").Value = 114518.6790515398
$ws.Range("E14").Value = 0.07520047540406788
$ws.Range("F14").Value = 0.1622905046453666
$ws.Range("G14").Value = -0.7872352754547525
$ws.Range("H14").Value = 11.27394741799627

$ws.Range("D16").Value = 114509.2964787599
$ws.Range("E16").Value = 0.06858969862743376
$ws.Range("F16").Value = 0.1590006229866709
$ws.Range("G16").Value = -1.010227271948733
$ws.Range("H16").Value = 10.2435996038116

$ws.Range("D17").Value = 114435.440610379
$ws.Range("E17").Value = 0.1459610303108773
$ws.Range("F17").Value = 0.1579052838648766
$ws.Range("G17").Value = -0.8365346370316579
$ws.Range("H17").Value = 5.637328271024248

$ws.Range("D20").Value = 115322.7329197569
$ws.Range("E20").Value = 0.04421666879072099
$ws.Range("F20").Value = 0.151938087846641
$ws.Range("G20").Value = -0.7790271005885648
$ws.Range("H20").Value = 7.477312722808811
